# Typo and broken link fixes
#
# The only user-facing content edit in this revision is a wording fix on
# slide 5 (the "Setting Up Free Body Diagrams" slide): the bullet reading
# "Assume Coefficient static and  of kinetic friction = .15 for all
# surfaces " is reworded to "Assume Coefficient static and  kinetic
# friction of .15 for all surfaces ".
#
# (The revision also shows a cached datetimeFigureOut value in the notes
# master and internal co-authoring/changesInfo + customXml bookkeeping
# churn, but those are artifacts PowerPoint itself rewrites on save/
# re-package rather than edits reachable through the Shape/Slide object
# model, so they are not something this script performs.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -like "Assume Coefficient static*") {
                $target = $shp
            }
        }
    }
}

if ($target -ne $null) {
    $target.TextFrame.TextRange.Text = "Assume Coefficient static and  kinetic friction of .15 for all surfaces "
}
